# Apply the parameter-comparison edits to the base DE/DG SW-upper
# parameter workbook: update the C13:D14 parameter values from 2 to 1,
# and move the active selection to C14 (as it was left after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "2" parameter values to "1" for rows 13 and 14 (columns C and D)
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 1

# Leave the selection on C14, matching the saved cursor position
$ws.Range("C14").Select()
